$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text so values like trailing-zero decimals
# and thousand-grouped numbers (e.g. 27.070.29) are preserved verbatim,
# matching the pre-existing text formatting of this column.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.070.29'
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').Value = '1.827.12'
$ws.Range('E3').Value = '  -0.73%  '
$ws.Range('E4').Value = '  -0.71%  '
$ws.Range('D5').Value = '311.48'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('D7').Value = '0.4254'
$ws.Range('E7').Value = '  -1.09%  '
$ws.Range('D8').Value = '0.3671'
$ws.Range('E8').Value = '  -1.52%  '
$ws.Range('D9').Value = '0.07234'
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('D10').Value = '0.8443'
$ws.Range('E10').Value = '  -3.00%  '
$ws.Range('D11').Value = '20.73'
$ws.Range('E11').Value = '  -3.00%  '
$ws.Range('D12').Value = '1.829.68'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').Value = '6.673'
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').Value = '5.293'
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = '0.07036'
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('D16').Value = '89.66'
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').Value = '0.000008752'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').Value = '14.89'
$ws.Range('E20').Value = '  -3.00%  '
$ws.Range('D21').Value = '27.146.40'
$ws.Range('E21').Value = '  -1.61%  '
$ws.Range('D22').Value = '5.134'
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('D24').Value = '2.053.00'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').Value = '1.983'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('D26').Value = '151.43'
$ws.Range('E26').Value = '  -1.99%  '
$ws.Range('D27').Value = '2.263'
$ws.Range('E27').Value = '  +5.19%  '
$ws.Range('D28').Value = '18.20'
$ws.Range('E28').Value = '  -1.67%  '
$ws.Range('D29').Value = '5.244'
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').Value = '116.78'
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('D31').Value = '0.08741'
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('D33').Value = '0.7374'
$ws.Range('E33').Value = '  -4.50%  '
$ws.Range('D34').Value = '4.430'
$ws.Range('E34').Value = '  -1.77%  '
$ws.Range('D35').Value = '2.901'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('D38').Value = '0.01944'
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('D39').Value = '0.05228'
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').Value = '7.308'
$ws.Range('E40').Value = '  +2.64%  '
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').Value = '0.5109'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '8.568'
$ws.Range('E44').Value = '  -1.77%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '10.54'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '1.977'
$ws.Range('E46').Value = '  +7.37%  '
$ws.Range('D47').Value = '0.4733'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = '105.54'
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('D51').Value = '1.653'
$ws.Range('E51').Value = '  -1.54%  '
